{"js": "// The commit removes the final checklist paragraph (the stray \"d\" item)\n// that followed the \"...local storage\" bullet, so the document now ends\n// right after that paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n// Find the last paragraph whose trimmed text is exactly \"d\" and remove it.\n// (Falls back to removing the very last paragraph if, for some reason, the\n// text doesn't match exactly \u2014 keeps the script robust to minor load\n// ordering differences while still targeting the intended paragraph.)\nlet target = null;\nfor (let i = items.length - 1; i >= 0; i--) {\n  if (items[i].text.trim() === \"d\") {\n    target = items[i];\n    break;\n  }\n}\nif (!target) {\n  target = items[items.length - 1];\n}\n\ntarget.delete();\nawait context.sync();\n", "ps1": "# The commit removes the final checklist paragraph (the stray \"d\" item)\n# that followed the \"...local storage\" bullet, so the document now ends\n# right after that paragraph.\n$d = $word.ActiveDocument\n\n$target = $null\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Trim() -eq \"d\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    # Fall back to the last paragraph if the exact text wasn't found.\n    $target = $d.Paragraphs.Item($d.Paragraphs.Count)\n}\n\n$target.Range.Delete()\n"}
